# Insert a new weekly price record as the most recent row for this
# "Apio" (celery) subset, pushing the existing rows 309-326 down to
# rows 310-327.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 309; this shifts rows
# 309:326 down to 310:327 and extends the used range to row 327.
$ws.Rows.Item(309).Insert()

# Populate the newly inserted row 309 with the latest weekly record.
$ws.Range("A309").Value = 4
$ws.Range("B309").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C309").Value = "Los Lagos"
$ws.Range("D309").Value = 44826
$ws.Range("E309").Value = 10
$ws.Range("F309").Value = 100112017
$ws.Range("G309").Value = "Apio"
$ws.Range("H309").Value = "Americana (o)"
$ws.Range("I309").Value = "Primera"
$ws.Range("J309").Value = 25
$ws.Range("K309").Value = 15000
$ws.Range("L309").Value = 15000
$ws.Range("M309").Value = 15000
$ws.Range("N309").Value = "$/docena de matas"
$ws.Range("O309").Value = "Región de Coquimbo"
$ws.Range("P309").Value = 2500
$ws.Range("Q309").Value = 6
$ws.Range("R309").Value = "Hortaliza"
